$d = $word.ActiveDocument

# Locate the first paragraph (the one with the intro sentence) and trim
# everything after the sentence "...to the same bookmark." (the trailing
# spacer run, the orange "<---" marker and the M2Doc version-mismatch
# message that were only there to flag a template/runtime version
# mismatch in the test fixture).
$p1 = $d.Paragraphs(1).Range

$rng = $d.Content
$rng.Start = $p1.Start
$rng.End = $p1.End

# Find the end of the sentence we want to keep.
$found = $rng.Find.Execute("same bookmark.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $cutStart = $rng.End
    $cutEnd = $p1.End - 1   # exclude the paragraph mark itself

    if ($cutEnd -gt $cutStart) {
        $toDelete = $d.Range($cutStart, $cutEnd)
        $toDelete.Text = ""
    }
}
